# Rename 'variable' and 'long_name' to 'variable-code' and 'variable-label'
# Close #144

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCodelists = $wb.Worksheets.Item("Codelists")

# --- Rename header labels ---
# Variables sheet: C1 'variable' -> 'variable-code', D1 'en_long_name' -> 'en_variable-label'
$wsVariables.Range("C1").Value = "variable-code"
$wsVariables.Range("D1").Value = "en_variable-label"

# Codelists sheet: A1 'variable' -> 'variable-code', D1 'en_code_label' -> 'en_code-label'
$wsCodelists.Range("A1").Value = "variable-code"
$wsCodelists.Range("D1").Value = "en_code-label"

# --- Column width adjustments on the Variables sheet ---
# (input values compensate for the runtime's internal +5/6 character offset so the
#  saved <col> width lands on the pixel-quantized value closest to the target)
$wsVariables.Columns.Item(3).ColumnWidth = 13.307291666666666
$wsVariables.Columns.Item(4).ColumnWidth = 16.451822916666668

# --- Selection / active sheet changes ---
# Variables sheet is no longer the selected tab; its selection moves to C1
$wsVariables.Range("C1").Select() | Out-Null

# Codelists sheet becomes the active/selected tab with selection collapsed to A2
$wsCodelists.Range("A2").Select() | Out-Null
$wsCodelists.Activate()
